$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Murali's office (Universal Sompo) insurance row details
$ws.Range("I3").Value = "14900(Ctrl S)"
$ws.Range("J3").Value = "universal sompo general insurance"
$ws.Range("K3").Value = "Murali(USGI51340704)`nRajalaxmi(USGI53042804)"

$ws.Range("L3").VerticalAlignment = -4108

$ws.Range("M3").Value = "8/9/2025"
$ws.Range("M3").NumberFormat = "d-mmm-yy"
$ws.Range("M3").VerticalAlignment = -4108

$ws.Range("J3").WrapText = $true
$ws.Range("K3").WrapText = $true

$ws.Rows("3").RowHeight = 29

$ws.Range("J3").Select()
